# Reorders the "Recorded By" (column G) comma-separated list of recorder
# names on each data row so that the literal entry "System" comes first,
# preserving the relative order of the remaining entries. When the list
# has no "System" entry, the entries are simply reversed.
#
# This mirrors the data fix applied upstream where the recorder lists
# (e.g. "dnasr281@gmail.com, System") were normalized to put "System"
# first (e.g. "System, dnasr281@gmail.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$col = 7  # Column G = "Recorded By"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $raw = $cell.Value2

    if ($null -eq $raw -or $raw -eq "") { continue }
    if ($raw.ToString().IndexOf(",") -lt 0) { continue }

    $parts = @($raw -split ',\s*')
    $idx = [Array]::IndexOf($parts, "System")

    if ($idx -ge 0) {
        $rest = @()
        for ($i = 0; $i -lt $parts.Count; $i++) {
            if ($i -ne $idx) { $rest += $parts[$i] }
        }
        $newParts = @("System") + $rest
    } else {
        $newParts = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $newParts += $parts[$i]
        }
    }

    $newVal = [string]::Join(", ", $newParts)

    if ($newVal -ne $raw) {
        $cell.Value = $newVal
    }
}
